$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Final Project" -> "Final Project " + "Report" (two runs; engine will
#    likely normalise to one run of identical formatting, which is fine).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Final Project", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Final Project Report", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Delete the "Phase 1:" heading paragraph entirely (bold one-liner).
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(6).Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 3) Replace the old Q1 paragraph text with the new intro paragraph.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(6)
$rng = $p.Range
$rng.End = $rng.End - 1
$rng.Text = "As many Texans are reeling from the impact of a major, statewide utility outage in February 2021, many conversations have been started to help utility providers and customers alike to ensure that this catastrophe does not occur again. One of the more prominent pushes from politicians and policy makers entails joining the nationwide power grid, instead of relying on a " + [char]8220 + "in house" + [char]8221 + " solution administered by ERCOT (Electric Reliability Council of Texas). This would be a impactful decision for the state of Texas, and all of living there, and this BI dashboard will attempt to help policy makers analyze and make the correct decision. "

# ---------------------------------------------------------------------------
# 4) Delete the Q2 paragraph, empty bold paragraph, "Phase 2, pt 1" heading,
#    the hyperlink paragraph, Username/Password/Database paragraphs — all of
#    them are collapsed into the "Tables..." paragraph which gets rewritten.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(8).Range.Delete() | Out-Null   # "2: What is the inflation..."
$d.Paragraphs.Item(8).Range.Delete() | Out-Null   # empty bold paragraph
$d.Paragraphs.Item(8).Range.Delete() | Out-Null   # "Phase 2, pt 1"
$d.Paragraphs.Item(8).Range.Delete() | Out-Null   # "Visit <hyperlink>"
$d.Paragraphs.Item(8).Range.Delete() | Out-Null   # "Username: root"
$d.Paragraphs.Item(8).Range.Delete() | Out-Null   # "Password: 1!vpnCwhite"
$d.Paragraphs.Item(8).Range.Delete() | Out-Null   # "Database: finalProject"

# Remove the now-orphaned hyperlink relationship / field, if any remain.
for ($i = $d.Hyperlinks.Count; $i -ge 1; $i--) {
    $d.Hyperlinks.Item($i).Delete() | Out-Null
}

$p = $d.Paragraphs.Item(8)
$rng = $p.Range
$rng.End = $rng.End - 1
$rng.Text = "These datasets are created by the US Department of Energy (a federal agency) to analyze and understand how different utility operators in the United States adjust their pricing and change operations based on environmental and social factors. While for the scope of this project, we only examine utilities in 2 states, Texas and Montana, given more time and resources, this BI dashboard could be expanded to include an analysis of all 50 states and give a bigger picture than what we see in our limited dashboard. Even in our limited dashboard, we are still able to see a trend in inflation rate, between two polar opposite states."

# ---------------------------------------------------------------------------
# 5) Delete "Phase 2, pt 2" heading, rewrite the "Given the time
#    constraints..." paragraph with the new content.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(10).Range.Delete() | Out-Null   # "Phase 2, pt 2"

$p = $d.Paragraphs.Item(10)
$rng = $p.Range
$rng.End = $rng.End - 1
$rng.Text = "First, we examined what the average residential utility rates are for two US states. This information will help new and upcoming providers, as well as existing to set rates that are applicable to everyone, and average and fair. Second, we examined the change, or inflation, between two different years, 2011 and 2019. To meet the requirements of the project, I used 2 data sources and created two database tables to store each. I first found the data describing utility providers and their characteristics from 2019. Given the time constraints, I was unable to identify a secondary dynamic data source that would be fetched via an API. To compensate for this, I identified a secondary data source from the Department of Energy in a near identical format to my original data which allowed me to make comparisons between the 2011 and 2019 years. While a real time dynamic source (such as real time rates from a utility in Texas, etc) might have provided a comparison between present day and 2019, I believe a comparison of 2011 and 2019 provides us with just as much insight. With proper development and implementation of data from additional states and years, this dashboard could become even more powerful."

# ---------------------------------------------------------------------------
# 6) Delete "Phase 3:" heading, rewrite "See attached pandas/jupyter file."
#    paragraph, then delete the trailing empty / "Project Report:" / empty
#    bold paragraphs.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(12).Range.Delete() | Out-Null   # "Phase 3:"

$p = $d.Paragraphs.Item(12)
$rng = $p.Range
$rng.End = $rng.End - 1
$rng.Text = "To share the information that is necessary, we first used Pandas to find the data for each state, and then return the average residential rate for each calendar year. This allows us to analyze the trends across different states while keeping in mind the particular challenges that a utility provider might have to overcome to provide quality service in each of these locations. Second, we use Pandas and Numpy to plot a bar graph which analyzes the difference in average rates between the two years. Finally, we use Pandas to calculate the precent change, or the inflation rate for the cost of electricity over the course of 8 years. The result is a targeted dashboard, backed by rich data to support the claims made."

$d.Paragraphs.Item(13).Range.Delete() | Out-Null   # empty paragraph
$d.Paragraphs.Item(13).Range.Delete() | Out-Null   # "Project Report:"
$d.Paragraphs.Item(13).Range.Delete() | Out-Null   # trailing empty bold paragraph

# ---------------------------------------------------------------------------
# 7) Remove now-unused character styles (Hyperlink / Unresolved Mention).
# ---------------------------------------------------------------------------
foreach ($styleName in @("Unresolved Mention", "Hyperlink")) {
    try {
        $style = $d.Styles.Item($styleName)
        if ($style -ne $null) {
            $style.Delete()
        }
    } catch {
    }
}

Write-Output "done"
